# Updated symbol list on Sat Jan 21 03:46:38 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '302.53' }
    @{ Cell = 'E2'; Value = '2.48%' }
    @{ Cell = 'D3'; Value = '35.09' }
    @{ Cell = 'E3'; Value = '13.27%' }
    @{ Cell = 'D4'; Value = '5.159' }
    @{ Cell = 'E4'; Value = '5.07%' }
    @{ Cell = 'D5'; Value = '0.07763' }
    @{ Cell = 'E5'; Value = '4.29%' }
    @{ Cell = 'D6'; Value = '2.340' }
    @{ Cell = 'E6'; Value = '7.67%' }
    @{ Cell = 'D7'; Value = '8.037' }
    @{ Cell = 'E7'; Value = '3.83%' }
    @{ Cell = 'B8'; Value = 'MXToken' }
    @{ Cell = 'C8'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' }
    @{ Cell = 'D8'; Value = '0.9311' }
    @{ Cell = 'E8'; Value = '1.65%' }
    @{ Cell = 'B9'; Value = 'LiechtensteinCryptoassetsExchange' }
    @{ Cell = 'C9'; Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx' }
    @{ Cell = 'D9'; Value = '0.09954' }
    @{ Cell = 'E9'; Value = '11.50%' }
    @{ Cell = 'B10'; Value = 'WazirX' }
    @{ Cell = 'C10'; Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx' }
    @{ Cell = 'D10'; Value = '0.1792' }
    @{ Cell = 'E10'; Value = '4.62%' }
    @{ Cell = 'B11'; Value = 'MandalaExchangeToken' }
    @{ Cell = 'C11'; Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx' }
    @{ Cell = 'D11'; Value = '0.08578' }
    @{ Cell = 'E11'; Value = '3.10%' }
    @{ Cell = 'B12'; Value = 'BitrueCoin' }
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr' }
    @{ Cell = 'D12'; Value = '0.03316' }
    @{ Cell = 'E12'; Value = '6.63%' }
    @{ Cell = 'B13'; Value = 'BitMartToken' }
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx' }
    @{ Cell = 'D13'; Value = '0.09917' }
    @{ Cell = 'E13'; Value = '-1.65%' }
    @{ Cell = 'B14'; Value = 'BitForexToken' }
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf' }
    @{ Cell = 'D14'; Value = '0.001502' }
    @{ Cell = 'E14'; Value = '-0.70%' }
    @{ Cell = 'B15'; Value = 'TigerCash' }
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch' }
    @{ Cell = 'D15'; Value = '0.005780' }
    @{ Cell = 'E15'; Value = '0.52%' }
    @{ Cell = 'B16'; Value = 'LEO' }
    @{ Cell = 'C16'; Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo' }
    @{ Cell = 'D16'; Value = '3.461' }
    @{ Cell = 'E16'; Value = '-1.23%' }
    @{ Cell = 'B17'; Value = 'GateToken' }
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt' }
    @{ Cell = 'D17'; Value = '3.945' }
    @{ Cell = 'E17'; Value = '5.35%' }
    @{ Cell = 'D18'; Value = '2.138' }
    @{ Cell = 'E18'; Value = '2.85%' }
    @{ Cell = 'D19'; Value = '0.3367' }
    @{ Cell = 'D20'; Value = '0.1334' }
    @{ Cell = 'E20'; Value = '2.68%' }
    @{ Cell = 'D21'; Value = '4.296' }
    @{ Cell = 'E21'; Value = '7.70%' }
    @{ Cell = 'E23'; Value = '-0.53%' }
    @{ Cell = 'E24'; Value = '0.13%' }
    @{ Cell = 'E25'; Value = '-5.24%' }
    @{ Cell = 'D26'; Value = '0.0001301' }
    @{ Cell = 'E26'; Value = '0.07%' }
    @{ Cell = 'E27'; Value = '-0.12%' }
    @{ Cell = 'D39'; Value = '0.01796' }
    @{ Cell = 'E39'; Value = '11.42%' }
    @{ Cell = 'D40'; Value = '0.04795' }
    @{ Cell = 'E40'; Value = '6.84%' }
    @{ Cell = 'D41'; Value = '0.007809' }
    @{ Cell = 'E41'; Value = '7.06%' }
    @{ Cell = 'E42'; Value = '6.03%' }
    @{ Cell = 'D43'; Value = '0.006838' }
    @{ Cell = 'E43'; Value = '-23.72%' }
    @{ Cell = 'D44'; Value = '0.002072' }
    @{ Cell = 'E44'; Value = '5.55%' }
    @{ Cell = 'D45'; Value = '0.009444' }
    @{ Cell = 'E45'; Value = '9.58%' }
    @{ Cell = 'D46'; Value = '0.00006111' }
    @{ Cell = 'E46'; Value = '0.31%' }
    @{ Cell = 'D47'; Value = '0.00000000750' }
    @{ Cell = 'E47'; Value = '-0.09%' }
    @{ Cell = 'D48'; Value = '3.006' }
    @{ Cell = 'E48'; Value = '34.79%' }
    @{ Cell = 'D49'; Value = '0.002001' }
    @{ Cell = 'E49'; Value = '-0.08%' }
    @{ Cell = 'D50'; Value = '0.00002101' }
    @{ Cell = 'E50'; Value = '-0.09%' }
    @{ Cell = 'D51'; Value = '0.0002001' }
    @{ Cell = 'E51'; Value = '-0.09%' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.ClearFormats()
}
